# Increase the gap between events dates
# Allow new events to be added between those events.

$wb = $excel.ActiveWorkbook

# --- confirmations (sheet1) -------------------------------------------------
$ws1 = $wb.Worksheets.Item("confirmations")
$ws1.Range("A3").Value = 44576
$ws1.Range("A4").Value = 44592
$ws1.Range("A5").Value = 44593
$ws1.Range("A6").Value = 44607
$ws1.Range("A7").Value = 44620
$ws1.Range("A8").Value = 44620
$ws1.Range("A9").Value = 44621
$ws1.Range("A10").Value = 44635

# --- trades (sheet2) ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("trades")
$ws2.Range("A3").Value = 44576
$ws2.Range("A4").Value = 44576
$ws2.Range("A5").Value = 44592
$ws2.Range("A6").Value = 44592
$ws2.Range("A7").Value = 44593
$ws2.Range("A8").Value = 44593
$ws2.Range("A9").Value = 44607
$ws2.Range("A10").Value = 44620
$ws2.Range("A11").Value = 44620
$ws2.Range("A12").Value = 44620
$ws2.Range("A13").Value = 44621
$ws2.Range("A14").Value = 44635
$ws2.Activate()
$ws2.Range("A2:A14").Select()

# --- subscriptions (sheet3) --------------------------------------------------
$ws3 = $wb.Worksheets.Item("subscriptions")
$ws3.Range("F2").Value = 44571
$ws3.Range("G2").Value = 44572
$ws3.Range("K2").Value = 44576
$ws3.Range("A3").Value = 44652
$ws3.Range("E3").Value = 44652
$ws3.Range("F3").Value = 44661
$ws3.Range("G3").Value = 44662

# --- splits (sheet4) ----------------------------------------------------------
$ws4 = $wb.Worksheets.Item("splits")
$ws4.Range("A2").Value = 44607
$ws4.Range("A3").Value = 44651
$ws4.Activate()
$ws4.Range("A3").Select()

# --- mergers (sheet5) ----------------------------------------------------------
$ws5 = $wb.Worksheets.Item("mergers")
$ws5.Range("A2").Value = 44651
$ws5.Activate()
$ws5.Range("A2").Select()

# --- spinoffs (sheet6) ----------------------------------------------------------
$ws6 = $wb.Worksheets.Item("spinoffs")
$ws6.Range("A2").Value = 44652

# --- restore the active sheet/tab to "confirmations" --------------------------
$ws1.Activate()
